$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.996.23"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "2.257.07"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'319.01"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'101.93"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").Value = "'0.578"
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.558"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'37.29"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "'0.0834"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "'7.70"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("D14").Value = "2.593.86"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "'0.861"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "'14.27"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").Value = "2.253.09"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "43.868.95"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").Value = "'13.50"
$ws.Range("E19").Value = "  -6.86%  "
$ws.Range("D20").Value = "0.0₃0987"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "'6.56"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "'65.81"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "'236.06"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").Value = "'2.13"
$ws.Range("E25").Value = "  -2.62%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'10.18"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'37.48"
$ws.Range("E28").Value = "  +2.95%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.16"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").Value = "'6.28"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").Value = "'20.25"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").Value = "'158.09"
$ws.Range("E32").Value = "  +3.78%  "
$ws.Range("D33").Value = "'0.0855"
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("D34").Value = "'2.71"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").Value = "'0.117"
$ws.Range("E35").Value = "  +11.98%  "
$ws.Range("D36").Value = "'3.10"
$ws.Range("E36").Value = "  -3.97%  "
$ws.Range("D37").Value = "'1.96"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'16.15"
$ws.Range("E39").Value = "  +16.50%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "'3.75"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").Value = "'4.23"
$ws.Range("E41").Value = "  -5.30%  "
$ws.Range("D42").Value = "'0.0317"
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "1.801.70"
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("D45").Value = "'0.200"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").Value = "'76.01"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'82.92"
$ws.Range("E47").Value = "  -4.35%  "
$ws.Range("D48").Value = "'5.23"
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.70"
$ws.Range("E49").Value = "  +6.63%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "'58.82"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'104.00"
$ws.Range("E51").Value = "  +0.56%  "
